# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp
# - Swap the Chequia/Noruega rows (country order change) and refresh their stats
# - Refresh Estados Unidos (row 4) and Brasil (row 11) stats

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 01:34"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1346771
$ws.Range("C4").Value = 24986
$ws.Range("D4").Value = 237022
$ws.Range("E4").Value = 1029722
$ws.Range("F4").Value = 16817
$ws.Range("G4").Value = 1412
$ws.Range("H4").Value = 80027

# Brasil (row 11)
$ws.Range("D11").Value = 61685
$ws.Range("E11").Value = 83627

# Row 49 becomes Noruega (was Chequia), with refreshed totals
$ws.Range("A49").Value = "Noruega"
$ws.Range("B49").Value = 8099
$ws.Range("C49").Value = 29
$ws.Range("D49").Value = 32
$ws.Range("E49").Value = 7848
$ws.Range("F49").Value = 24
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 219

# Row 50 becomes Chequia (was Noruega), with refreshed totals
$ws.Range("A50").Value = "Chequia"
$ws.Range("B50").Value = 8095
$ws.Range("C50").Value = 18
$ws.Range("D50").Value = 4446
$ws.Range("E50").Value = 3373
$ws.Range("F50").Value = 47
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 276
